$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated capital structure database values for cost_equity (X), roe_cost_equity (Y),
# cost_capital (AB) and roic_cost_capital (AC) columns, rows 2-7.

$ws.Range("X2").Value = 0.1755391088970272
$ws.Range("Y2").Value = 0.2146017361733953
$ws.Range("AB2").Value = 0.1249859114804322
$ws.Range("AC2").Value = -0.1863109341003176

$ws.Range("X3").Value = 0.2509871020863186
$ws.Range("Y3").Value = 0.1391537429841039
$ws.Range("AB3").Value = 0.1320898544138489
$ws.Range("AC3").Value = -0.1476800548592832

$ws.Range("X4").Value = 0.1695145060287309
$ws.Range("Y4").Value = 0.5081963373447631
$ws.Range("AB4").Value = 0.1243419199509988
$ws.Range("AC4").Value = -0.1525554993840904

$ws.Range("X5").Value = 0.1755391088970272
$ws.Range("Y5").Value = 0.1958894625315443
$ws.Range("AB5").Value = 0.1249859114804322
$ws.Range("AC5").Value = -0.1863109341003176

$ws.Range("X6").Value = 0.2658515462999966
$ws.Range("Y6").Value = 0.1394289817528087
$ws.Range("AB6").Value = 0.1306574189839833
$ws.Range("AC6").Value = -0.392323884677942

$ws.Range("X7").Value = 0.1153115367055389
$ws.Range("Y7").Value = -0.8733760528345712
$ws.Range("AB7").Value = 0.1143435684093837
$ws.Range("AC7").Value = -0.8387251938510798
